$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# This workbook holds one worksheet per OpenDSS "class", each sheet sharing
# the same 39-column (A:AM) header template (Id_<Class>, bus1, bus2, ...).
# The commit adds four brand-new (empty, header-only) class sheets:
#   - WindGen          inserted immediately before "GenDispatcher"
#   - GICLine          inserted immediately before "GICTransformer"
#   - FMonitor         appended right after "Sensor"
#   - Generic5         appended right after "FMonitor"
# plus one unrelated data fix on the "Capacitor" sheet (cell H2).
#
# New sheets are created with Worksheet.Copy so they inherit the exact same
# sheetPr/pageMargins/column styling (s="1" header style) as their template,
# then are renamed and have their A1 header ("Id_...") corrected.
#
# NOTE: worksheet object references become stale (re-point to a different
# sheet by position) once another sheet is inserted/copied, so every sheet
# reference is re-fetched by name with Worksheets.Item(...) right before use.
# NOTE: this PowerShell shim does not bind named (-Param value) arguments to
# function params reliably, so positional params are used everywhere.
# ---------------------------------------------------------------------------

function Add-ClassSheet {
    param($TemplateName, $NewName, $InsertBefore)

    $template = $wb.Worksheets.Item($TemplateName)
    if ($InsertBefore) {
        # Copy(Before, After) -> duplicate lands right before the template;
        # the duplicate is the one named "<TemplateName> (2)".
        $template.Copy($template, $null)
    } else {
        # duplicate lands right after the template
        $template.Copy($null, $template)
    }

    $dup = $wb.Worksheets.Item($TemplateName + " (2)")
    $dup.Name = $NewName

    $dup = $wb.Worksheets.Item($NewName)
    $dup.Range("A1").Value = "Id_" + $NewName
}

# 1) WindGen, right before GenDispatcher
Add-ClassSheet "GenDispatcher" "WindGen" $true

# 2) GICLine, right before GICTransformer
Add-ClassSheet "GICTransformer" "GICLine" $true

# 3) FMonitor, right after Sensor
Add-ClassSheet "Sensor" "FMonitor" $false

# 4) Generic5, right after FMonitor
Add-ClassSheet "FMonitor" "Generic5" $false

# 5) Unrelated data fix: Capacitor!H2 literal text value
$cap = $wb.Worksheets.Item("Capacitor")
$cap.Range("H2").Value = "('0 |1.87639338887875E-310 2.25607960651843E-308 |1.15409172934718E-305 0 1.15409061474694E-305',)"

Write-Output "Sheets now:"
$wb.Worksheets | ForEach-Object { Write-Output $_.Name }
